# Applies the textual edits described by the commit diff to the resume.
# Uses Find/Execute replacements scoped narrowly (and with MatchCase) so
# that only the intended occurrence of each string is touched.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText, [bool]$matchCase = $true) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $matchCase, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> [$findText]"
    }
    return $ok
}

# 1. "Desenvolvedor Back-end" -> "Desenvolvedor back-end" (lowercase the "Back-")
Replace-Text " Back-" " back-"

# 2. Fix "Habilidaes" typo -> "Habilidades"
Replace-Text "Habilidaes" "Habilidades"

# 3. Spring boot proficiency "3" years -> "1" year. ": 3" also appears
#    earlier for the unrelated "Java: 3" skill, so locate "Spring boot"
#    first and only search-and-replace in the remainder of the document
#    from that point on (this also keeps the bold "Spring boot" run and
#    the grey ": 1" run separate, matching the original split).
$springRng = $d.Content
$springFound = $springRng.Find.Execute("Spring boot", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
if ($springFound) {
    $afterSpring = $d.Range($springRng.End, $d.Content.End)
    $ok = $afterSpring.Find.Execute(": 3", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, ": 1", 1)
    if (-not $ok) {
        Write-Host "WARNING: not found -> [: 3 after Spring boot]"
    }
} else {
    Write-Host "WARNING: not found -> [Spring boot]"
}

# 4. "Coordenação grupos de trabalho." -> "Coordenação de grupos de trabalho."
Replace-Text "Coordenação grupos de trabalho." "Coordenação de grupos de trabalho."

# 5. "Disp. trabalho remoto" -> "Disp. p/ trabalho remoto"
Replace-Text "Disp. trabalho remoto" "Disp. p/ trabalho remoto"

# 6. Lowercase "Inspeção" / "Operação" mid-sentence
Replace-Text "Inspeção e organização do ambiente de operação (5s), Operação de equipamentos" `
             "inspeção e organização do ambiente de operação (5s), operação de equipamentos"

# 7. Explicitly mark the section as portrait orientation (already portrait
#    by dimensions, this just writes the attribute Word adds on resave).
$d.PageSetup.Orientation = 0
